$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "28÷5=" "85÷7="
Replace-Text "25÷9=" "19÷6="
Replace-Text "55÷6=" "50÷3="
Replace-Text "24÷9=" "20÷8="
Replace-Text "97÷5=" "78÷5="
Replace-Text "38÷9=" "14÷2="
Replace-Text "51÷4=" "45÷8="
Replace-Text "93÷9=" "90÷9="
Replace-Text "36÷2=" "10÷9="
Replace-Text "64÷2=" "79÷9="
Replace-Text "24÷3=" "96÷5="
Replace-Text "19÷4=" "30÷2="
Replace-Text "95÷4=" "11÷7="
Replace-Text "64÷9=" "24÷5="
Replace-Text "79÷6=" "25÷9="
Replace-Text "86÷2=" "21÷5="
Replace-Text "73÷4=" "93÷5="
Replace-Text "79÷8=" "84÷3="
Replace-Text "43÷2=" "92÷5="
Replace-Text "91÷5=" "68÷3="
Replace-Text "69÷5=" "70÷4="
Replace-Text "26÷7=" "78÷5="
Replace-Text "90÷6=" "26÷5="
Replace-Text "30÷4=" "22÷9="
Replace-Text "42÷7=" "98÷2="
